$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Market Cap" column (C) with refreshed figures for rows 2-26
$ws.Range("C2").Value = 734283209566.5956
$ws.Range("C3").Value = 245227573040.3853
$ws.Range("C4").Value = 39201564034.85326
$ws.Range("C5").Value = 33310656530.94685
$ws.Range("C6").Value = 24399622165.73147
$ws.Range("C7").Value = 13665400136.30362
$ws.Range("C8").Value = 11304976145.07183
$ws.Range("C9").Value = 9056483331.556614
$ws.Range("C10").Value = 8260159516.992143
$ws.Range("C11").Value = 8084726899.33145
$ws.Range("C12").Value = 7632285366.632638
$ws.Range("C13").Value = 7618079701.282753
$ws.Range("C14").Value = 6698546429.708444
$ws.Range("C15").Value = 6132020529.331023
$ws.Range("C16").Value = 5186461309.489183
$ws.Range("C17").Value = 5048987793.845691
$ws.Range("C18").Value = 4465611717.579144
$ws.Range("C19").Value = 3658175713.09136
$ws.Range("C20").Value = 3461649966.767
$ws.Range("C21").Value = 3434662981.505833
$ws.Range("C22").Value = 3363821194.031478
$ws.Range("C23").Value = 3012735565.923972
$ws.Range("C24").Value = 2993198303.274961
$ws.Range("C25").Value = 2851150359.793557
$ws.Range("C26").Value = 2432964843.350379
